$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row34
$ws.Cells.Item(34, 8).Value = 5000
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 14).ClearContents()

# ALC!row36
$ws.Cells.Item(36, 8).Value = 5000
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 14).ClearContents()

# ALC!row40
$ws.Cells.Item(40, 8).Value = 3980
$ws.Cells.Item(40, 9).Value = 2950
$ws.Cells.Item(40, 11).Value = 2950
$ws.Cells.Item(40, 13).Value = -2775

# ALC!row51
$ws.Cells.Item(51, 8).Value = 9597.1
$ws.Cells.Item(51, 9).Value = 7794.4
$ws.Cells.Item(51, 10).Value = 11399.8
$ws.Cells.Item(51, 11).Value = 7794.4
$ws.Cells.Item(51, 12).Value = 11399.8
$ws.Cells.Item(51, 13).Value = -7310.4
$ws.Cells.Item(51, 14).Value = -12367.8

$ws = $wb.Worksheets.Item("ARM")
# ARM!row5
$ws.Cells.Item(5, 8).Value = 303.75
$ws.Cells.Item(5, 9).Value = 303.75
$ws.Cells.Item(5, 11).Value = 303.75
$ws.Cells.Item(5, 13).Value = -191.75

# ARM!row32
$ws.Cells.Item(32, 8).Value = 5032.625
$ws.Cells.Item(32, 9).Value = 5032.625
$ws.Cells.Item(32, 11).Value = 5032.625
$ws.Cells.Item(32, 13).Value = -4745.625

# ARM!row41
$ws.Cells.Item(41, 8).Value = 676.25
$ws.Cells.Item(41, 9).Value = 676.25
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 11).Value = 676.25
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 13).ClearContents()
$ws.Cells.Item(41, 14).Value = -262.25

$ws = $wb.Worksheets.Item("BSM")
# BSM!row4
$ws.Cells.Item(4, 8).Value = 303.75
$ws.Cells.Item(4, 9).Value = 303.75
$ws.Cells.Item(4, 11).Value = 303.75
$ws.Cells.Item(4, 13).Value = -188.75

# BSM!row97
$ws.Cells.Item(97, 8).Value = 9439.6
$ws.Cells.Item(97, 9).Value = 9439.6
$ws.Cells.Item(97, 11).Value = 9439.6
$ws.Cells.Item(97, 13).Value = -8448.6

$ws = $wb.Worksheets.Item("CRP")
# CRP!row15
$ws.Cells.Item(15, 8).Value = 11663.333
$ws.Cells.Item(15, 10).Value = 11663.333
$ws.Cells.Item(15, 12).Value = 11663.333
$ws.Cells.Item(15, 14).Value = -12003.333

# CRP!row58
$ws.Cells.Item(58, 8).Value = 971.625
$ws.Cells.Item(58, 9).Value = 953.2857
$ws.Cells.Item(58, 10).Value = 1100
$ws.Cells.Item(58, 11).Value = 953.2857
$ws.Cells.Item(58, 12).Value = 1100
$ws.Cells.Item(58, 13).Value = -750.2857
$ws.Cells.Item(58, 14).Value = -1506

# CRP!row134
$ws.Cells.Item(134, 8).Value = 3199.8
$ws.Cells.Item(134, 9).Value = 3000
$ws.Cells.Item(134, 10).Value = 3499.5
$ws.Cells.Item(134, 11).Value = 9000
$ws.Cells.Item(134, 12).Value = 10498.5
$ws.Cells.Item(134, 13).Value = -6465
$ws.Cells.Item(134, 14).Value = -15568.5

# CRP!row136
$ws.Cells.Item(136, 8).Value = 971.625
$ws.Cells.Item(136, 9).Value = 953.2857
$ws.Cells.Item(136, 10).Value = 1100
$ws.Cells.Item(136, 11).Value = 2859.8571
$ws.Cells.Item(136, 12).Value = 3300
$ws.Cells.Item(136, 13).Value = -309.8571000000002
$ws.Cells.Item(136, 14).Value = -8400

$ws = $wb.Worksheets.Item("CUL")
# CUL!row3
$ws.Cells.Item(3, 8).Value = 9998
$ws.Cells.Item(3, 9).Value = 9997
$ws.Cells.Item(3, 11).Value = 29991
$ws.Cells.Item(3, 13).Value = -29879

# CUL!row10
$ws.Cells.Item(10, 8).Value = 271.72726
$ws.Cells.Item(10, 9).Value = 20.31579
$ws.Cells.Item(10, 10).Value = 1864
$ws.Cells.Item(10, 11).Value = 60.94737
$ws.Cells.Item(10, 12).Value = 5592
$ws.Cells.Item(10, 13).Value = 78.05262999999999
$ws.Cells.Item(10, 14).Value = -5870

# CUL!row13
$ws.Cells.Item(13, 8).Value = 434.5625
$ws.Cells.Item(13, 9).Value = 90
$ws.Cells.Item(13, 11).Value = 270
$ws.Cells.Item(13, 13).Value = -102

# CUL!row25
$ws.Cells.Item(25, 8).Value = 178.18182
$ws.Cells.Item(25, 9).Value = 123.125
$ws.Cells.Item(25, 10).Value = 325
$ws.Cells.Item(25, 11).Value = 369.375
$ws.Cells.Item(25, 12).Value = 975
$ws.Cells.Item(25, 13).Value = -200.375
$ws.Cells.Item(25, 14).Value = -1313

# CUL!row30
$ws.Cells.Item(30, 8).Value = 178.18182
$ws.Cells.Item(30, 9).Value = 123.125
$ws.Cells.Item(30, 10).Value = 325
$ws.Cells.Item(30, 11).Value = 369.375
$ws.Cells.Item(30, 12).Value = 975
$ws.Cells.Item(30, 13).Value = -267.375
$ws.Cells.Item(30, 14).Value = -1179

# CUL!row37
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 12).Value = 0
$ws.Cells.Item(37, 14).ClearContents()

# CUL!row119
$ws.Cells.Item(119, 8).Value = 1676.3334
$ws.Cells.Item(119, 9).Value = 1676.3334
$ws.Cells.Item(119, 11).Value = 5029.0002
$ws.Cells.Item(119, 13).Value = -191.0002000000004

# CUL!row133
$ws.Cells.Item(133, 8).Value = 2500
$ws.Cells.Item(133, 9).Value = 2500
$ws.Cells.Item(133, 11).Value = 7500
$ws.Cells.Item(133, 13).Value = -2440

# CUL!row134
$ws.Cells.Item(134, 8).Value = 8333.333000000001
$ws.Cells.Item(134, 9).Value = 5000
$ws.Cells.Item(134, 11).Value = 15000
$ws.Cells.Item(134, 13).Value = -9930

# CUL!row138
$ws.Cells.Item(138, 8).Value = 1500000
$ws.Cells.Item(138, 9).Value = 1500000
$ws.Cells.Item(138, 11).Value = 4500000
$ws.Cells.Item(138, 13).Value = -4494860

$ws = $wb.Worksheets.Item("GSM")
# GSM!row132
$ws.Cells.Item(132, 8).Value = 941.6667
$ws.Cells.Item(132, 9).Value = 912.5
$ws.Cells.Item(132, 10).Value = 1000
$ws.Cells.Item(132, 11).Value = 2737.5
$ws.Cells.Item(132, 12).Value = 3000
$ws.Cells.Item(132, 13).Value = -207.5
$ws.Cells.Item(132, 14).Value = -8060

$ws = $wb.Worksheets.Item("LTW")
# LTW!row7
$ws.Cells.Item(7, 8).Value = 10155.648
$ws.Cells.Item(7, 9).Value = 9994.161
$ws.Cells.Item(7, 11).Value = 9994.161
$ws.Cells.Item(7, 13).Value = -9882.161

# LTW!row40
$ws.Cells.Item(40, 8).Value = 1400
$ws.Cells.Item(40, 9).Value = 1400
$ws.Cells.Item(40, 11).Value = 1400
$ws.Cells.Item(40, 13).Value = -1264

# LTW!row62
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 14).ClearContents()

# LTW!row65
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 14).ClearContents()

# LTW!row93
$ws.Cells.Item(93, 8).Value = 33335920
$ws.Cells.Item(93, 9).Value = 66669148
$ws.Cells.Item(93, 11).Value = 66669148
$ws.Cells.Item(93, 13).Value = -66667900

# LTW!row126
$ws.Cells.Item(126, 8).Value = 10155.648
$ws.Cells.Item(126, 9).Value = 9994.161
$ws.Cells.Item(126, 11).Value = 29982.483
$ws.Cells.Item(126, 13).Value = -27512.483

$ws = $wb.Worksheets.Item("WVR")
# WVR!row7
$ws.Cells.Item(7, 8).Value = 799
$ws.Cells.Item(7, 9).Value = 799
$ws.Cells.Item(7, 11).Value = 799
$ws.Cells.Item(7, 13).Value = -686

# WVR!row97
$ws.Cells.Item(97, 8).Value = 46000
$ws.Cells.Item(97, 10).Value = 46000
$ws.Cells.Item(97, 12).Value = 46000
$ws.Cells.Item(97, 14).Value = -47982

# WVR!row126
$ws.Cells.Item(126, 8).Value = 4282.923
$ws.Cells.Item(126, 9).Value = 3387.8
$ws.Cells.Item(126, 11).Value = 10163.4
$ws.Cells.Item(126, 13).Value = -7693.400000000001
